{"js": "\n// ---------------------------------------------------------------------\n// 1) \"Fabio work differently ... reloading.\" was split into two runs by\n//    a stray _GoBack bookmark. Re-merge it into a single run and drop\n//    the bookmark.\n// ---------------------------------------------------------------------\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet fabioPara = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Fabio work differently\") !== -1) {\n    fabioPara = p;\n    break;\n  }\n}\n\nif (fabioPara) {\n  const fabioInnerXml = \"<w:p w:rsidR=\\\"00C772AB\\\" w:rsidRDefault=\\\"00C772AB\\\" w:rsidP=\\\"00C772AB\\\"><w:pPr><w:pStyle w:val=\\\"ListParagraph\\\"/><w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"2\\\"/></w:numPr></w:pPr><w:r><w:t>Fabio work differently since it updates it routing table directly from the data stored in Consul as soon as there is a change and without restart and reloading.</w:t></w:r></w:p>\";\n  const fabioOoxml =\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n    fabioInnerXml +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>';\n  fabioPara.insertOoxml(fabioOoxml, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2) Append the new \"Jaeger Tracing\" section (intro text + a bulleted\n//    list of advantages, ending with two blank trailer paragraphs) right\n//    after the last paragraph of the document\n//    (\"...-/order and Fabio will do the test.\").\n// ---------------------------------------------------------------------\nconst allParagraphs = body.paragraphs;\nallParagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = allParagraphs.items[allParagraphs.items.length - 1];\nconst placeholder = lastParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\n// A trailing fully-empty paragraph (<w:p/>) placed at the very end of the\n// body tends to get absorbed into the body's closing paragraph mark when\n// imported via insertOoxml, so the chunk below ends with a throwaway\n// sentinel paragraph that is deleted afterwards, leaving the real empty\n// paragraph as the new last paragraph in the body.\nconst jaegerInnerXml = \"<w:p><w:r><w:t xml:space=\\\"preserve\\\">Jaeger Tracing </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\\\"preserve\\\">Jaeger is a CNCF distributed tracing system released by Uber that enables debugging, monitoring, and analysis of your services, and is based on </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>OpenTracing</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> API. It uses distributed context propagation, which essentially is the basis of distributed tracing, to assign metadata to request as they propagate through your system. List some advantages of using Jaeger:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"ListParagraph\\\"/><w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"2\\\"/></w:numPr></w:pPr><w:r><w:t>Performance and latency analysis.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"ListParagraph\\\"/><w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"2\\\"/></w:numPr></w:pPr><w:r><w:t>Service dependency analysis; you can view a DAG of your system in the UI</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"ListParagraph\\\"/><w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"2\\\"/></w:numPr></w:pPr><w:r><w:t>Logs associated with each span</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"ListParagraph\\\"/><w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"2\\\"/></w:numPr></w:pPr><w:r><w:t>Organization of logs into calling hierarchy</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"ListParagraph\\\"/><w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"2\\\"/></w:numPr></w:pPr><w:r><w:t>Cost attribution</w:t></w:r><w:bookmarkStart w:id=\\\"0\\\" w:name=\\\"_GoBack\\\"/><w:bookmarkEnd w:id=\\\"0\\\"/></w:p><w:p><w:pPr><w:pStyle w:val=\\\"ListParagraph\\\"/><w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"2\\\"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Useful UI and libraries in Go, Node, C#, and others.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"ListParagraph\\\"/><w:numPr><w:ilvl w:val=\\\"0\\\"/><w:numId w:val=\\\"2\\\"/></w:numPr></w:pPr><w:r><w:t xml:space=\\\"preserve\\\">Different levels of sampling can be configured to reduce load: </w:t></w:r><w:proofErr w:type=\\\"gramStart\\\"/><w:r><w:t>constant ,</w:t></w:r><w:proofErr w:type=\\\"gramEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> probabilistic, rate limiting, and remote.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\\\"ListParagraph\\\"/></w:pPr></w:p><w:p/><w:p><w:r><w:t>ZZZ_SENTINEL_DELETE_ME_ZZZ</w:t></w:r></w:p>\";\nconst jaegerOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n  jaegerInnerXml +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\nplaceholder.insertOoxml(jaegerOoxml, Word.InsertLocation.replace);\nawait context.sync();\n\nconst finalParagraphs = body.paragraphs;\nfinalParagraphs.load(\"items/text\");\nawait context.sync();\n\nlet sentinel = null;\nfor (const p of finalParagraphs.items) {\n  if (p.text.indexOf(\"ZZZ_SENTINEL_DELETE_ME_ZZZ\") !== -1) {\n    sentinel = p;\n    break;\n  }\n}\nif (sentinel) {\n  sentinel.delete();\n  await context.sync();\n}\n", "ps1": "# -----------------------------------------------------------------------\n# 1) \"Fabio work differently ... reloading.\" is currently split across two\n#    runs by a stray \"_GoBack\" bookmark sitting between them. Delete that\n#    bookmark and re-run Find/Replace over the sentence so Word merges it\n#    back into a single run.\n# -----------------------------------------------------------------------\n$d = $word.ActiveDocument\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$sentence = \"Fabio work differently since it updates it routing table directly from the data stored in Consul as soon as there is a change and without restart and reloading.\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute($sentence, $false, $false, $false, $false, $false, $true, 1, $false, $sentence, 2) | Out-Null\n\n# -----------------------------------------------------------------------\n# 2) Append the new \"Jaeger Tracing\" section (intro paragraphs plus a\n#    bulleted list of advantages, ending with two blank trailer\n#    paragraphs) right after the existing last paragraph of the document\n#    (\"...-/order and Fabio will do the test.\").\n# -----------------------------------------------------------------------\n$endRange = $d.Content\n$endRange.Collapse(0) | Out-Null\n\n$jaegerXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t xml:space=\"preserve\">Jaeger Tracing </w:t></w:r></w:p><w:p><w:r><w:t xml:space=\"preserve\">Jaeger is a CNCF distributed tracing system released by Uber that enables debugging, monitoring, and analysis of your services, and is based on </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>OpenTracing</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> API. It uses distributed context propagation, which essentially is the basis of distributed tracing, to assign metadata to request as they propagate through your system. List some advantages of using Jaeger:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr><w:r><w:t>Performance and latency analysis.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr><w:r><w:t>Service dependency analysis; you can view a DAG of your system in the UI</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr><w:r><w:t>Logs associated with each span</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr><w:r><w:t>Organization of logs into calling hierarchy</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr><w:r><w:t>Cost attribution</w:t></w:r><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Useful UI and libraries in Go, Node, C#, and others.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr></w:pPr><w:r><w:t xml:space=\"preserve\">Different levels of sampling can be configured to reduce load: </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>constant ,</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> probabilistic, rate limiting, and remote.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"ListParagraph\"/></w:pPr></w:p><w:p/>'\n$endRange.InsertXML($jaegerXml) | Out-Null\n"}
